# "RF classify (read data header bug)"
#
# The sheet used to keep a duplicate "raw data" block in columns H:M
# (a second copy of the A:F header plus the raw clock-time readings),
# and A2:B7 were formulas that subtracted a "Video start time:" anchor
# (stored in A9/B9) from those raw H:M readings to get elapsed times.
# That raw/anchor scaffolding was a read-data-header bug: the real
# output is just the elapsed-time values in A:F. Clean it up:
#   1. freeze A2:B7 to their computed values (so nothing breaks once
#      their precedent columns disappear),
#   2. delete the duplicate raw-data columns H:M,
#   3. delete the now-unused "Video start time:" row (row 9), which
#      shifts the lone constant that used to live in row 11 up to
#      row 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Convert A2:B7 from formulas (=H2-B$9, =I2-B$9, ...) to the static
#    values they already evaluate to, before their precedents go away.
$rng = $ws.Range("A2:B7")
$rng.Value2 = $rng.Value2

# 2) Remove the duplicate raw-data header/readings block in H1:M7.
$ws.Range("H1:M7").EntireColumn.Delete() | Out-Null

# 3) Remove the "Video start time:" row; everything below (the B11
#    constant) shifts up by one row, landing on row 10.
$ws.Rows(9).Delete() | Out-Null

# Leave the selection where the saved workbook shows it.
$ws.Range("H1:O1048576").Select() | Out-Null
